$wb = $excel.ActiveWorkbook

# "Write data" sheet (sheet5.xml): update CorpDashGO value in B4
$wsWrite = $wb.Worksheets.Item("Write data")
$wsWrite.Range("B4").Value = "Offered - 0; Onboarded - 0; TotalMoneySaved - Coming Soon"

# "Login Page" sheet (sheet1.xml): update Admin Email value in B4
$wsLogin = $wb.Worksheets.Item("Login Page")
$wsLogin.Range("B4").Value = "automate@workstreets.com"

# Restore active sheet / selection to match final workbook state (Login Page, B5 selected)
$wsLogin.Activate()
$wsLogin.Range("B5").Select()
